$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.164600729942322
$ws.Range("B1").Value = 2.141740322113037
$ws.Range("C1").Value = 3.347212791442871
$ws.Range("D1").Value = 3.644415378570557
$ws.Range("E1").Value = 1.178421974182129
